# Update the NEW_HIGHEST_ALL sheet: Samsung NEW pricelist refresh
# (Z Fold/Flip/S25 family price changes, new S25 Edge rows, removal of
#  Galaxy A17 4G/5G and Galaxy A26 5G rows, per commit message.)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW_HIGHEST_ALL")

$data = @(
    ,@("Samsung", "Galaxy Z Fold 7 5G", "256GB", 1630)
    ,@("Samsung", "Galaxy Z Fold 7 5G", "512GB", 1780)
    ,@("Samsung", "Galaxy Z Fold 7 5G", "1TB", 1860)
    ,@("Samsung", "Galaxy Z Flip 7 5G", "256GB", 900)
    ,@("Samsung", "Galaxy Z Flip 7 5G", "512GB", 950)
    ,@("Samsung", "Galaxy Z Flip 7 FE 5G", "128GB", 620)
    ,@("Samsung", "Galaxy Z Flip 7 FE 5G", "256GB", 770)
    ,@("Samsung", "Galaxy S25 5G", "256GB", 850)
    ,@("Samsung", "Galaxy S25 5G", "512GB", 950)
    ,@("Samsung", "Galaxy S25+ 5G", "256GB", 920)
    ,@("Samsung", "Galaxy S25+ 5G", "512GB", 1120)
    ,@("Samsung", "Galaxy S25 Ultra 5G", "256GB", 1020)
    ,@("Samsung", "Galaxy S25 Ultra 5G", "512GB", 1200)
    ,@("Samsung", "Galaxy S25 Ultra 5G", "1TB", 1350)
    ,@("Samsung", "Galaxy S25 Edge 5G", "256GB", 770)
    ,@("Samsung", "Galaxy S25 Edge 5G", "512GB", 870)
    ,@("Samsung", "Galaxy S25 FE 5G", "128GB", 520)
    ,@("Samsung", "Galaxy S25 FE 5G", "256GB", 620)
    ,@("Samsung", "Galaxy S25 FE 5G", "512GB", 670)
    ,@("Samsung", "Galaxy A36 5G", "8/256GB", 340)
    ,@("Samsung", "Galaxy A56 5G", "8/256GB", 380)
    ,@("Samsung", "Galaxy A56 5G", "12/256GB", 420)
    ,@("Samsung", "Galaxy Tab A11+ 128 WiFi", "Base", 260)
    ,@("Samsung", "Galaxy Tab A11+ 128 5G", "Base", 320)
    ,@("Samsung", "Galaxy Watch 8 40mm Bluetooth", "Base", 280)
    ,@("Samsung", "Galaxy Watch 8 44mm Bluetooth", "Base", 300)
    ,@("Samsung", "Galaxy Watch 8 Classic 46mm Bluetooth", "Base", 360)
    ,@("Samsung", "Galaxy Watch Ultra 47mm (2025)", "Base", 480)
    ,@("Samsung", "Galaxy Buds 3", "Base", 50)
    ,@("Samsung", "Galaxy Buds 3 Pro", "Base", 130)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# The old sheet had 31 data rows (through row 32); the refreshed list has
# only 30 data rows (through row 31), so clear out the now-unused last row.
$lastOldRow = $startRow + $data.Count
$ws.Range($ws.Cells.Item($lastOldRow, 1), $ws.Cells.Item($lastOldRow, 4)).Clear()
